# Add 2022-Q3 data
# -----------------------------------------------------------------------
# The workbook has a "总计" (overview) sheet followed by one sheet per
# quarter (most-recent-first). This change inserts a brand new "2022-Q3"
# quarter sheet right after "总计" (pushing every older quarter sheet
# down by one tab position, but leaving their data untouched), fills it
# with that quarter's fund-holding data, and adds a matching summary row
# to the "总计" sheet.

$wb = $excel.ActiveWorkbook

$zongji = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q3" worksheet right after "总计".
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $zongji)
$q3.Name = "2022-Q3"

# Header row (matches the other quarter sheets).
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"
$q3.Range("B1:H1").Font.Bold = $true
$q3.Range("B1:H1").Borders.LineStyle = 1
$q3.Range("B1:H1").HorizontalAlignment = -4108
$q3.Range("B1:H1").VerticalAlignment = -4160

# Fund-code / numeric-text columns must stay text (fund codes like
# "166005" and the numeric-looking figures are stored as text in every
# quarter sheet), so force Text format before writing them.
$q3.Range("B2:B9").NumberFormat = "@"
$q3.Range("D2:G9").NumberFormat = "@"

$q3Data = @(
    @(0, "166005", "中欧价值发现混合 -A",              "26.62", "93.73", "7.63", "2.0311", 2),
    @(1, "001810", "中欧潜力价值灵活配置混合A",          "19.07", "93.66", "7.80", "1.4875", 2),
    @(2, "004232", "中欧价值发现混合 -C",              "8.18",  "93.73", "7.63", "0.6241", 2),
    @(3, "166024", "中欧恒利三年定期开放混合",           "3.99",  "98.45", "6.95", "0.2773", 2),
    @(4, "005764", "中欧潜力价值灵活配置混合C",          "2.01",  "93.66", "7.80", "0.1568", 2),
    @(5, "166020", "中欧成长优选回报灵活配置混合A",       "2.38",  "93.70", "5.76", "0.1371", 2),
    @(6, "001891", "中欧成长优选回报灵活配置混合E",       "0.74",  "93.70", "5.76", "0.0426", 2),
    @(7, "001882", "中欧价值发现混合 -E",              "0.43",  "93.73", "7.63", "0.0328", 2)
)

$r = 2
foreach ($row in $q3Data) {
    $q3.Cells.Item($r, 1).Value = $row[0]
    $q3.Cells.Item($r, 2).Value = $row[1]
    $q3.Cells.Item($r, 3).Value = $row[2]
    $q3.Cells.Item($r, 4).Value = $row[3]
    $q3.Cells.Item($r, 5).Value = $row[4]
    $q3.Cells.Item($r, 6).Value = $row[5]
    $q3.Cells.Item($r, 7).Value = $row[6]
    $q3.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

# Column A (the running index) uses the same bold/bordered style as the
# other quarter sheets.
$q3.Range("A2:A9").Font.Bold = $true
$q3.Range("A2:A9").Borders.LineStyle = 1
$q3.Range("A2:A9").HorizontalAlignment = -4108
$q3.Range("A2:A9").VerticalAlignment = -4160

$q3.Range("A1").Select()

# ---------------------------------------------------------------------
# 2. Add the 2022-Q3 summary row to "总计", shifting the older rows down.
# ---------------------------------------------------------------------
$zongji.Rows.Item(2).Insert()
$zongji.Range("A2:D2").ClearFormats()

$zongji.Range("A3").Copy()
$zongji.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$zongji.Range("A2").Value = 0
$zongji.Range("B2").Value = "2022-Q3"
$zongji.Range("C2").Value = 8
$zongji.Range("D2").Value = 4.79

# Renumber the running index (column A) of the pre-existing rows, which
# all shifted down by one.
for ($row = 3; $row -le 9; $row++) {
    $zongji.Cells.Item($row, 1).Value = $row - 2
}

$zongji.Range("A1").Select()
